$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 50: new "Undergraduate Enrollment" label (mirrors B35's formatting) ---
$ws.Cells.Item(35, 2).Copy($ws.Cells.Item(50, 2))

# --- New "Evanston, IL" / Northwestern section (rows 58-64) ---

# Row 58: Northwestern University enrollment figure
$ws.Cells.Item(58, 1).Value = "Northwestern University"
$ws.Cells.Item(58, 2).Value = 23161

# Row 59: Hilton Garden Inn Chicago North Shore
$ws.Cells.Item(59, 1).Value = "Hilton Garden Inn Chicago North Shore"
$ws.Cells.Item(59, 2).Value = 178
$ws.Cells.Item(59, 4).Value = "https://www.choosechicago.com/listing/hilton-garden-inn-chicago-north-shore-evanston/"
$ws.Hyperlinks.Add($ws.Cells.Item(59, 4), "https://www.choosechicago.com/listing/hilton-garden-inn-chicago-north-shore-evanston/", "", "", "")

# Row 60: Hilton Orrington Evanston
$ws.Cells.Item(60, 1).Value = "Hilton Orrington Evanston"
$ws.Cells.Item(60, 2).Value = 269
$ws.Cells.Item(60, 4).Value = "https://www.hilton.com/en/hotels/ordoehf-hilton-orrington-evanston/events/"
$ws.Hyperlinks.Add($ws.Cells.Item(60, 4), "https://www.hilton.com/en/hotels/ordoehf-hilton-orrington-evanston/events/", "", "", "")

# Row 61: Graduate Evanston
$ws.Cells.Item(61, 1).Value = "Graduate Evanston"
$ws.Cells.Item(61, 2).Value = 119
$ws.Cells.Item(61, 4).Value = "https://www.choosechicago.com/listing/graduate-evanston/"
$ws.Hyperlinks.Add($ws.Cells.Item(61, 4), "https://www.choosechicago.com/listing/graduate-evanston/", "", "", "")

# Row 62: Hyatt House Chicago Evanston (hyperlink carries a #HotelName fragment as "location")
$ws.Cells.Item(62, 1).Value = "Hyatt House Chicago Evanston"
$ws.Cells.Item(62, 2).Value = 114
$ws.Cells.Item(62, 4).Value = "https://www.hotelplanner.com/Hotels/223659/Reservations-Hyatt-House-Chicago-Evanston-Evanston-1515-Chicago-Ave-60201#HotelName"
$ws.Hyperlinks.Add($ws.Cells.Item(62, 4), "https://www.hotelplanner.com/Hotels/223659/Reservations-Hyatt-House-Chicago-Evanston-Evanston-1515-Chicago-Ave-60201", "HotelName", "", "")

# Row 63: GEM Museum Suites
$ws.Cells.Item(63, 1).Value = "GEM Museum Suites"
$ws.Cells.Item(63, 2).Value = 71
$ws.Cells.Item(63, 4).Value = "https://www.qantas.com/hotels/properties/1126680-the-gem-museum-suites"
$ws.Hyperlinks.Add($ws.Cells.Item(63, 4), "https://www.qantas.com/hotels/properties/1126680-the-gem-museum-suites", "", "", "")

# Row 64: TOTAL = SUM(B59:B63)
$ws.Cells.Item(64, 1).Value = "TOTAL"
$ws.Cells.Item(64, 2).Formula = "=SUM(B59:B63)"

# --- View state: scroll down and select A69 ---
$ws.Range("A69").Select()
$excel.ActiveWindow.ScrollRow = 58
